# PowerShell Excel COM-interop script
# Applies price / volume(1h) updates to the cryptos worksheet,
# including the row-swaps for Hedera/Maker (rows 38-39) and Cosmos/SuiNetwork (rows 50-51).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# --- Per-row Price (D) / Volume 1h (E) updates ---
Set-TextValue "D2" "64.909.73"
Set-TextValue "E2" "  +2.17%  "
Set-TextValue "D3" "3.468.39"
Set-TextValue "E3" "  +1.68%  "
Set-TextValue "E4" "  +0.04%  "
Set-TextValue "D5" "575.82"
Set-TextValue "E5" "  +1.18%  "
Set-TextValue "D6" "161.60"
Set-TextValue "E6" "  +3.19%  "
Set-TextValue "D7" "0.615"
Set-TextValue "E7" "  +13.02%  "
Set-TextValue "D9" "3.472.18"
Set-TextValue "E9" "  +1.81%  "
Set-TextValue "D10" "7.22"
Set-TextValue "E10" "  -2.14%  "
Set-TextValue "D11" "0.124"
Set-TextValue "E11" "  +2.10%  "
Set-TextValue "D12" "0.445"
Set-TextValue "E12" "  +3.55%  "
Set-TextValue "D13" "4.073.42"
Set-TextValue "E13" "  +1.87%  "
Set-TextValue "E14" "  +0.46%  "
Set-TextValue "D15" "0.0000191"
Set-TextValue "E15" "  -0.79%  "
Set-TextValue "D16" "28.29"
Set-TextValue "E16" "  +4.00%  "
Set-TextValue "D17" "65.003.00"
Set-TextValue "E17" "  +2.25%  "
Set-TextValue "D18" "3.457.73"
Set-TextValue "E18" "  +1.79%  "
Set-TextValue "E19" "  +3.26%  "
Set-TextValue "D20" "14.34"
Set-TextValue "E20" "  +2.03%  "
Set-TextValue "D21" "381.15"
Set-TextValue "E21" "  +1.04%  "
Set-TextValue "D22" "8.14"
Set-TextValue "E22" "  +0.60%  "
Set-TextValue "D23" "0.550"
Set-TextValue "E23" "  +4.22%  "
Set-TextValue "D24" "72.50"
Set-TextValue "E24" "  +1.08%  "
Set-TextValue "E25" "  -0.01%  "
Set-TextValue "D26" "0.0000119"
Set-TextValue "E26" "  -0.80%  "
Set-TextValue "D27" "10.07"
Set-TextValue "E27" "  +7.24%  "
Set-TextValue "E28" "  +0.53%  "
Set-TextValue "E29" "  +0.10%  "
Set-TextValue "D30" "1.51"
Set-TextValue "E30" "  +10.89%  "
Set-TextValue "D31" "6.14"
Set-TextValue "E31" "  +1.21%  "
Set-TextValue "E32" "  +1.94%  "
Set-TextValue "D33" "23.59"
Set-TextValue "E33" "  +1.45%  "
Set-TextValue "D34" "7.13"
Set-TextValue "E34" "  +5.39%  "
Set-TextValue "E35" "  +12.80%  "
Set-TextValue "D36" "161.79"
Set-TextValue "E36" "  +1.95%  "
Set-TextValue "E37" "  +6.32%  "
Set-TextValue "D40" "26.62"
Set-TextValue "E40" "  -1.61%  "
Set-TextValue "D41" "6.72"
Set-TextValue "E41" "  +6.00%  "
Set-TextValue "D42" "4.55"
Set-TextValue "E42" "  +5.71%  "
Set-TextValue "D43" "0.0320"
Set-TextValue "E43" "  +1.48%  "
Set-TextValue "D44" "42.87"
Set-TextValue "E44" "  +2.28%  "
Set-TextValue "D45" "0.780"
Set-TextValue "E45" "  +2.51%  "
Set-TextValue "D46" "25.73"
Set-TextValue "E46" "  +11.02%  "
Set-TextValue "D47" "1.10"
Set-TextValue "E47" "  +3.83%  "
Set-TextValue "D48" "315.19"
Set-TextValue "E48" "  +7.01%  "
Set-TextValue "D49" "0.110"
Set-TextValue "E49" "  +7.50%  "

# --- Row 38/39 swap: Maker <-> Hedera (Hedera now ranks above Maker) ---
Set-TextValue "B38" "Hedera"
Set-TextValue "C38" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D38" "0.0777"
Set-TextValue "E38" "  +2.55%  "

Set-TextValue "B39" "Maker"
Set-TextValue "C39" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D39" "2.973.34"
Set-TextValue "E39" "  -0.20%  "

# --- Row 50/51 swap: SuiNetwork <-> Cosmos (Cosmos now ranks above SuiNetwork) ---
Set-TextValue "B50" "Cosmos"
Set-TextValue "C50" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D50" "6.63"
Set-TextValue "E50" "  +4.80%  "

Set-TextValue "B51" "SuiNetwork"
Set-TextValue "C51" "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue "D51" "0.868"
Set-TextValue "E51" "  +4.84%  "
